# Append a new job listing row (row 3) to the "ランサーズ" sheet and refresh
# the "取得日時" (fetched-at) timestamp on every existing data row to the
# latest run time: 2026-01-17 12:35:59.
#
# Net effect vs. the original sheet:
#   - row 2 keeps its data, only the timestamp changes
#   - a brand-new row is inserted as row 3
#   - the old rows 3-8 shift down to become rows 4-9 (timestamp refreshed)
#   - the hyperlinks in column F are rebuilt so they keep pointing at the
#     correct URL for each (possibly shifted) row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2026-01-17 12:35:59"

# Final desired contents for data rows 2..9, in sheet order.
# Columns: A=取得日時 B=タイトル C=カテゴリ D=価格 E=締切 F=URL G=優先度スコア H=スキル概要
$rows = @(
    @($timestamp, "【急募】airtableで社内業務管理システムを共に構築してくれる方", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5473383", 353, "🔥AI,Ai ◇管理"),
    @($timestamp, "【シンプル版】生成AIデジタル・コミュニティ制作の依頼", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469128", 310, "🔥AI,Ai"),
    @($timestamp, "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5217096", 243, "🔥API ◆ツール"),
    @($timestamp, "※急募:Next.jsによる業務アプリの開発(+Flutter)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5473147", 225, "🔥Next.js ◆開発 ◇アプリ"),
    @($timestamp, "※急募:Flutterによる業務アプリの開発(+next.js)", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5473146", 218, "🔥Next.js ◆開発 ◇アプリ"),
    @($timestamp, "【急募】Accessでの受発注管理・請求書発行システム開発", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5473234", 148, "◆開発,システム開発 ◇管理"),
    @($timestamp, "【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5473394", 135, "◆ツール,開発"),
    @($timestamp, "【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発", "システム開発", "1,000,000 円 ~ 3,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5473181", 135, "◆ツール,開発")
)

# Write every cell value first (this also naturally extends the sheet's
# dimension from H8 to H9).
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $vals = $rows[$i]
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
    $ws.Range("D$r").Value = $vals[3]
    $ws.Range("E$r").Value = $vals[4]
    $ws.Range("F$r").Value = $vals[5]
    $ws.Range("G$r").Value = $vals[6]
    $ws.Range("H$r").Value = $vals[7]
}

# Rebuild the column-F hyperlinks so they line up with the (shifted) URLs.
# Clearing hyperlinks via any single range removes every hyperlink on the
# sheet, so do it once up front and then re-add them all in final order.
$ws.Range("F2").Hyperlinks.Delete()

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $url = $rows[$i][5]
    $ws.Hyperlinks.Add($ws.Range("F$r"), $url)
}
